$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Month"
$ws.Range("C1").Value = "Accommodation"
$ws.Range("D1").Value = "Category"
$ws.Range("E1").Value = "Amount"
$ws.Range("F1").Value = "Description"

$ws.Range("E10").Select()
